# Updated symbol list on Fri Dec 23 04:22:16 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and a couple of "Volume(1h)" (column E)
# values in the crypto price table with the latest scraped figures.
# Column D values are stored as text (not numbers), so a leading
# apostrophe is used to force Excel to keep them as text instead of
# silently coercing the numeric-looking strings into Number cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'247.03"
$ws.Range("D3").Value  = "'22.02"
$ws.Range("D4").Value  = "'5.432"
$ws.Range("D6").Value  = "'3.405"
$ws.Range("D7").Value  = "'6.329"
$ws.Range("D8").Value  = "'0.8189"
$ws.Range("D9").Value  = "'0.9924"
$ws.Range("E9").Value  = "8FTXTokenFTTBestin24h"
$ws.Range("D10").Value = "'0.1431"
$ws.Range("D11").Value = "'0.07464"
$ws.Range("D12").Value = "'0.03144"
$ws.Range("D13").Value = "'0.02995"
$ws.Range("D14").Value = "'4.166"
$ws.Range("D15").Value = "'0.09412"
$ws.Range("D16").Value = "'0.001590"
$ws.Range("D17").Value = "'0.04825"
$ws.Range("D18").Value = "'0.0005852"
$ws.Range("D19").Value = "'0.006193"
$ws.Range("D21").Value = "'0.0009968"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.776"
$ws.Range("D24").Value = "'2.219"
$ws.Range("D25").Value = "'0.3230"
$ws.Range("D27").Value = "'0.0004000"
$ws.Range("D40").Value = "'0.03898"
$ws.Range("D41").Value = "'0.006463"
$ws.Range("D43").Value = "'0.002631"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.006718"
$ws.Range("D45").Value = "'0.00005594"
$ws.Range("D47").Value = "'0.3801"
$ws.Range("D49").Value = "'0.00002101"
